# Insert a new weekly record row into the Hortaliza / Zapallo sheet.
# A new data row is inserted at row 356, pushing the existing rows
# 356-370 down to 357-371 (dimension grows from A1:R370 to A1:R371).
# The new row duplicates the record that is now in row 357 (same
# market / region / product / quality) except for the date, the
# quality label and the "Volumen" (J) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(356).Insert()

$ws.Range("A356").Value = 5
$ws.Range("B356").Value = "Macroferia Regional de Talca"
$ws.Range("C356").Value = "Maule"
$ws.Range("D356").Value = 44939
$ws.Range("E356").Value = 7
$ws.Range("F356").Value = 100112045
$ws.Range("G356").Value = "Zapallo"
$ws.Range("H356").Value = "Camote"
$ws.Range("I356").Value = "1a nueva(o)"
$ws.Range("J356").Value = 900
$ws.Range("K356").Value = 450
$ws.Range("L356").Value = 450
$ws.Range("M356").Value = 450
$ws.Range("N356").Value = "$/kilo (volumen en unidades)"
$ws.Range("O356").Value = "Región del Maule"
$ws.Range("P356").Value = 450
$ws.Range("Q356").Value = 1
$ws.Range("R356").Value = "Hortaliza"
